# "changes ib ds and register"
#
# - tryEditorCode!A3: replace the old placeholder text with the start of a
#   pasted-in snippet (leading tab + "print('Hello").
# - Make "tryEditorCode" the active sheet/tab (previously "Register" was
#   active), with A3 selected there.

$wb = $excel.ActiveWorkbook

$wsEditor = $wb.Worksheets.Item("tryEditorCode")
$wsEditor.Range("A3").Value = "`tprint('Hello"

$wsEditor.Activate() | Out-Null
$wsEditor.Range("A3").Select() | Out-Null
